$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updated values
$ws.Range("B2").Value = 0.003994804209775715
$ws.Range("C2").Value = 0.0000006633126561350622
$ws.Range("D2").Value = 0.8054896365839992
$ws.Range("E2").Value = 645.3272768299601
$ws.Range("G2").Value = 646.1367619340665

# Row 3 updated values
$ws.Range("B3").Value = 3.230985683306322
$ws.Range("C3").Value = 1.667794583268128
$ws.Range("D3").Value = 0.1575252929769615
$ws.Range("E3").Value = 0.496779210170732
$ws.Range("G3").Value = 5.553084769722144
